$d = $word.ActiveDocument

# --- Locate the contact-info paragraph (the "Body Text" styled paragraph holding the
#     "{phone_number} | {email} | ..." placeholder line). ---
$p = $null
foreach ($cand in $d.Paragraphs) {
    if ($cand.Range.Text -like "*{phone_number}*{email}*") {
        $p = $cand
    }
}
if ($p -eq $null) {
    $p = $d.Paragraphs(2)
}

$pStart = $p.Range.Start
$pEnd = $p.Range.End

# --- Collapse "{phone_number} | {email} | {github} | {linkedin} | {portfolio} | {location}"
#     down to "{header_line}", reusing the existing "{" run that precedes "email" (and its
#     formatting) plus the run that held "email" (renamed to "header_line}"). ---

# 1) Remove everything before the "{" that opens "{email}" (i.e. the phone_number block and its
#    trailing " | ").
$find1 = $d.Range($pStart, $pEnd)
$find1.Find.Execute("{email}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$braceStart = $find1.Start

$lead = $d.Range($pStart, $braceStart)
$lead.Delete()

# 2) Remove everything after the closing "}" of "{email}" through the end of the paragraph text
#    (the " | {github} | {linkedin} | {portfolio} | {location}" tail).
$pEndNow = $p.Range.End
$find2 = $d.Range($p.Range.Start, $pEndNow)
$find2.Find.Execute("{email}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$emailCloseEnd = $find2.End

$trail = $d.Range($emailCloseEnd, $pEndNow - 1)
$trail.Delete()

# 3) Rename the placeholder from "email" to "header_line".
$find3 = $d.Range($p.Range.Start, $p.Range.End)
$find3.Find.Execute("email", $true, $false, $false, $false, $false, $true, 1, $false, "header_line", 2)

# 4) The leading "{" and "header_line}" now share identical run formatting, so the engine merges
#    them into a single run. Nudge a character-level property on the first character and revert
#    it so the run boundary between "{" and "header_line}" is preserved as two separate runs.
$firstChar = $d.Range($p.Range.Start, $p.Range.Start + 1)
$firstChar.Font.Bold = 1
$firstChar.Font.Bold = 0

# --- Paragraph spacing: before=60 (unchanged), after 40 -> 60, line 193/exact -> 300/auto ---
$pf = $p.Range.ParagraphFormat
$pf.LineSpacingRule = 5   # wdLineSpaceMultiple
$pf.LineSpacing = 15
$pf.SpaceAfter = 3
